$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.904.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.875.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7426"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3156"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.55%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07219"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.72%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08407"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7513"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.68%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.882.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.893.35"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.080"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.57"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007850"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9994"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.125.40"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.024"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1562"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.75%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.506"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.615"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.532"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.280"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05318"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.17%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7531"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9966"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.689"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01962"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.759"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4524"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.112.36"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.044"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.59"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8548"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.49"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.857"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.624"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.464"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.023.30"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.67%  "
